# Refresh coin prices / 1h volume change percentages pulled from coinranking.com
# (scheduled GitHub Actions data sync)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '33.933.92'
$ws.Range("E2").Value = '  -0.17%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '1.782.50'
$ws.Range("E3").Value = '  -0.21%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.18%  '
# Row 5: BNB
$ws.Range("D5").Value = '''226.08'
$ws.Range("E5").Value = '  +2.20%  '
# Row 6: XRP
$ws.Range("D6").Value = '''0.545'
$ws.Range("E6").Value = '  -1.31%  '
# Row 7: USDC
$ws.Range("E7").Value = '  +0.21%  '
# Row 8: Solana
$ws.Range("D8").Value = '''31.97'
$ws.Range("E8").Value = '  -1.09%  '
# Row 9: Cardano
$ws.Range("D9").Value = '''0.292'
$ws.Range("E9").Value = '  +3.23%  '
# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.0678'
$ws.Range("E10").Value = '  -4.56%  '
# Row 11: TRON
$ws.Range("E11").Value = '  +0.97%  '
# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '2.038.98'
$ws.Range("E12").Value = '  +0.05%  '
# Row 13: Chainlink
$ws.Range("D13").Value = '''11.15'
$ws.Range("E13").Value = '  +3.51%  '
# Row 14: WrappedEther
$ws.Range("D14").Value = '1.777.91'
$ws.Range("E14").Value = '  +0.10%  '
# Row 15: WrappedBTC
$ws.Range("D15").Value = '33.907.45'
$ws.Range("E15").Value = '  -0.16%  '
# Row 16: Polygon
$ws.Range("D16").Value = '''0.614'
$ws.Range("E16").Value = '  -1.59%  '
# Row 17: Polkadot
$ws.Range("E17").Value = '  +0.25%  '
# Row 18: Litecoin
$ws.Range("D18").Value = '''67.47'
$ws.Range("E18").Value = '  -0.32%  '
# Row 19: BitcoinCash
$ws.Range("D19").Value = '''241.55'
$ws.Range("E19").Value = '  -0.77%  '
# Row 20: ShibaInu
$ws.Range("D20").Value = '0.0₃0769'
$ws.Range("E20").Value = '  -1.55%  '
# Row 21: Dai
$ws.Range("E21").Value = '  -0.08%  '
# Row 22: Avalanche
$ws.Range("D22").Value = '''10.63'
$ws.Range("E22").Value = '  -1.59%  '
# Row 23: Uniswap
$ws.Range("D23").Value = '''4.07'
$ws.Range("E23").Value = '  +0.02%  '
# Row 24: Toncoin
$ws.Range("D24").Value = '''2.06'
$ws.Range("E24").Value = '  -2.06%  '
# Row 25: Monero
$ws.Range("D25").Value = '''161.71'
$ws.Range("E25").Value = '  +2.58%  '
# Row 26: Cosmos
$ws.Range("D26").Value = '''7.13'
$ws.Range("E26").Value = '  +1.48%  '
# Row 27: EthereumClassic
$ws.Range("D27").Value = '''16.15'
$ws.Range("E27").Value = '  -1.02%  '
# Row 28: Stellar
$ws.Range("D28").Value = '''0.112'
$ws.Range("E28").Value = '  +0.22%  '
# Row 29: BinanceUSD
$ws.Range("E29").Value = '  +0.32%  '
# Row 30: PancakeSwap
$ws.Range("D30").Value = '''1.23'
$ws.Range("E30").Value = '  +2.87%  '
# Row 31: Hedera
$ws.Range("E31").Value = '  -0.98%  '
# Row 32: Filecoin
$ws.Range("D32").Value = '''3.61'
$ws.Range("E32").Value = '  -1.17%  '
# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = '''3.55'
$ws.Range("E33").Value = '  +1.65%  '
# Row 34: LidoDAOToken
$ws.Range("E34").Value = '  +0.55%  '
# Row 35: Maker
$ws.Range("D35").Value = '1.392.38'
$ws.Range("E35").Value = '  -0.12%  '
# Row 36: ImmutableX
$ws.Range("D36").Value = '''0.644'
$ws.Range("E36").Value = '  +0.88%  '
# Row 37: TrustWalletToken
$ws.Range("E37").Value = '  -1.15%  '
# Row 38: VeChain
$ws.Range("E38").Value = '  +1.25%  '
# Row 39: RenderToken
$ws.Range("E39").Value = '  +8.50%  '
# Row 40: Aave
$ws.Range("D40").Value = '''79.60'
$ws.Range("E40").Value = '  +0.20%  '
# Row 41: HuobiToken
$ws.Range("E41").Value = '  +0.50%  '
# Row 42: ARBITRUM
$ws.Range("D42").Value = '''0.917'
$ws.Range("E42").Value = '  -0.32%  '
# Row 43: MXToken
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '''13.63'
$ws.Range("E43").Value = '  +14.27%  '
# Row 44: InjectiveProtocol
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = '''2.67'
$ws.Range("E44").Value = '  -1.56%  '
# Row 45: BabyDogeCoin
$ws.Range("D45").Value = '0.0₆0141'
$ws.Range("E45").Value = '  +11.59%  '
# Row 46: WEMIXToken
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").Value = '''0.0509'
$ws.Range("E46").Value = '  +3.40%  '
# Row 47: Kaspa
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '''1.08'
$ws.Range("E47").Value = '  +2.76%  '
# Row 48: FraxShare
$ws.Range("D48").Value = '''5.90'
$ws.Range("E48").Value = '  +0.66%  '
# Row 49: Quant
$ws.Range("D49").Value = '''107.05'
$ws.Range("E49").Value = '  +0.01%  '
# Row 50: RocketPoolETH
$ws.Range("D50").Value = '1.941.30'
$ws.Range("E50").Value = '  +0.02%  '
# Row 51: PaxDollar
$ws.Range("E51").Value = '  +0.19%  '
